# Updated cryptos list values (Price and Volume(1h) columns) per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.317.65"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.497.20"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.92"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.72"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.70"
$ws.Range("E9").Value = "  +7.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Value = "4.097.77"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "3.499.91"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "64.279.93"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.10"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.99"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.77"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.53"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.85"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.578"
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("D23").Value = "3.639.13"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.32"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.74"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.49"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.12"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("D34").Value = "3.527.93"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.28"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.41"
$ws.Range("E37").Value = "  +4.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "164.82"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0786"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.807"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.41"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.29"
$ws.Range("E46").Value = "  -4.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.64"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").Value = "2.436.69"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.79"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.920"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("E51").Value = "  -0.38%  "
